# Add a team "record" (Wins / Losses / Ties) to the roster sheet.
# The existing data occupies columns A:AC (rows 1-58, row 1 = headers).
# We extend it with three new columns: AD=Wins, AE=Losses, AF=Ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (AD1:AF1) -------------------------------------------
# Copy the formatting of the last existing header cell (AC1, style index 1:
# bold, bordered, centered) onto the new header cells before setting their
# text, so the new headers look consistent with the rest of the header row.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- New data columns (rows 2-58) -----------------------------------------
# Every player on the roster shares the same team record.
$ws.Range("AD2:AD58").Value = 88
$ws.Range("AE2:AE58").Value = 74
$ws.Range("AF2:AF58").Value = 0

$excel.CutCopyMode = 0
